$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column corrections (rows 20, 21, 28, 34) ---
# These were all "finished" before; now a couple of projects turn out to
# have no repo / couldn't be imported.
$ws.Range("C20").Value = "No Git Repo"
$ws.Range("C21").Value = "can't import"
$ws.Range("C28").Value = "can't import"
$ws.Range("C34").Value = "can't import"

# --- Newly assigned projects: rows 58-78 get owner + status ---
for ($r = 58; $r -le 78; $r++) {
    $ws.Cells.Item($r, 2).Value = "Keye Li"
    $ws.Cells.Item($r, 3).Value = "IN PROGRESS"
}

# --- Widen the status column slightly to fit the new text ---
$ws.Columns("C").ColumnWidth = 12.1

# --- Update the view: scroll position + selected cell ---
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
$ws.Range("N69").Select() | Out-Null
